# Add supplier / Group logic
# Updates the sample data rows (4-6) of the InventoryList table with new
# Id / Name / Supplier / Price / Threshold / Interval / Group values, and
# restores the selection left active in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InventoryList")

# --- Row 4 ---------------------------------------------------------------
$ws.Range("C4").Value = "Id01"
$ws.Range("D4").Value = "Tovar01"
$ws.Range("E4").Value = "Sup1"
$ws.Range("F4").Value = 40
$ws.Range("I4").Value = 32
$ws.Range("J4").Value = 20
$ws.Range("M4").Value = "Group2"

# --- Row 5 ---------------------------------------------------------------
$ws.Range("C5").Value = "Id02"
$ws.Range("D5").Value = "Tovar02"
$ws.Range("E5").Value = "Sup2"
$ws.Range("F5").Value = 22
$ws.Range("I5").Value = 32
$ws.Range("J5").Value = 14
$ws.Range("M5").Value = "Group1"

# --- Row 6 -----------------------------------------------------------------
$ws.Range("C6").Value = "Id03"
$ws.Range("D6").Value = "Tovar03"
$ws.Range("E6").Value = "Sup2"
$ws.Range("F6").Value = 43
$ws.Range("I6").Value = 33
$ws.Range("J6").Value = 55
$ws.Range("M6").Value = "Group2"

# --- Selection -------------------------------------------------------------
$ws.Activate()
$ws.Range("G6").Select() | Out-Null
